# ------------------------------------------------------------------
# Quick-save checkpoint edit: fill in the LR parse-table example with
# actual state/action data, extend the blank table body down to row
# 37, add a new "ITEM SET/STATE #" key entry, and restyle the header /
# highlighted rows (bold font, yellow/orange fills, thin borders).
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1) Header row (A1:H1) + the merged/label cell A2: bold font, no
#    underline, orange fill (already present), thin border all round.
# ---------------------------------------------------------------
$header = $ws.Range("A1:H1")
$header.Font.Bold = $true
$header.Font.Underline = $false
$header.Borders.LineStyle = 1
$header.Borders.Weight = 2
$header.Borders.Color = 0

$a2 = $ws.Range("A2")
$a2.Font.Bold = $true
$a2.Font.Underline = $false
$a2.Borders.LineStyle = 1
$a2.Borders.Weight = 2
$a2.Borders.Color = 0

# ---------------------------------------------------------------
# 2) Row 2 column headers (0 / 1 / * / + / $ / E / B): bold font,
#    yellow highlight fill, thin border all round.
# ---------------------------------------------------------------
$bc2 = $ws.Range("B2:C2")
$bc2.Font.Bold = $true
$bc2.Interior.Color = 65535
$bc2.Borders.LineStyle = 1
$bc2.Borders.Weight = 2
$bc2.Borders.Color = 0

$def2 = $ws.Range("D2:F2")
$def2.Font.Bold = $true
$def2.Interior.Color = 65535
$def2.Borders.LineStyle = 1
$def2.Borders.Weight = 2
$def2.Borders.Color = 0

$gh2 = $ws.Range("G2:H2")
$gh2.Font.Bold = $true
$gh2.Interior.Color = 65535
$gh2.Borders.LineStyle = 1
$gh2.Borders.Weight = 2
$gh2.Borders.Color = 0

# ---------------------------------------------------------------
# 3) Parse-table body values (states / actions / goto entries).
# ---------------------------------------------------------------
$ws.Range("A3").Value = "I0"
$ws.Range("B3").Value = "S5"
$ws.Range("C3").Value = "S6"
$ws.Range("D3").Value = "ERR"
$ws.Range("E3").Value = "ERR"
$ws.Range("F3").Value = "ERR"
$ws.Range("G3").Value = "S1"
$ws.Range("H3").Value = "S4"

$ws.Range("A4").Value = "I1"
$ws.Range("B4").Value = "ERR"
$ws.Range("C4").Value = "ERR"
$ws.Range("D4").Value = "S2"
$ws.Range("E4").Value = "S3"
$ws.Range("F4").Value = "ACC"
$ws.Range("G4").Value = "ERR"
$ws.Range("H4").Value = "ERR"

$ws.Range("A5").Value = "I2"
$ws.Range("B5").Value = "ERR"
$ws.Range("C5").Value = "ERR"
$ws.Range("D5").Value = "ERR"
$ws.Range("E5").Value = "ERR"
$ws.Range("F5").Value = "ERR"
$ws.Range("G5").Value = "ERR"
$ws.Range("H5").Value = "S7"

$ws.Range("A6").Value = "I3"
$ws.Range("B6").Value = "ERR"
$ws.Range("C6").Value = "ERR"
$ws.Range("D6").Value = "ERR"
$ws.Range("E6").Value = "ERR"
$ws.Range("F6").Value = "ERR"
$ws.Range("G6").Value = "ERR"
$ws.Range("H6").Value = "S8"

$ws.Range("A7").Value = "I4"

# ---------------------------------------------------------------
# 4) Alignment for the body cells (matches how the author filled
#    them in: center+vcenter for most of the first couple columns,
#    plain center for the rest of the grid / blank filler rows).
# ---------------------------------------------------------------
$vcenterCells = "A3,B3,C3,D3,E3,A4,B4,C4,A5,B5,C5,A6,B6,C6,A7"
foreach ($addr in $vcenterCells.Split(",")) {
    $ws.Range($addr).HorizontalAlignment = -4108
    $ws.Range($addr).VerticalAlignment = -4108
}

$centerOnlyCells = New-Object System.Collections.Generic.List[string]
foreach ($addr in @("F3","G3","H3")) { [void]$centerOnlyCells.Add($addr) }
foreach ($r in 4..37) {
    foreach ($col in @("D","E","F","G","H")) {
        [void]$centerOnlyCells.Add("$col$r")
    }
}
foreach ($addr in $centerOnlyCells) {
    $ws.Range($addr).HorizontalAlignment = -4108
}

# ---------------------------------------------------------------
# 5) New "ITEM SET/STATE #" entry in the table key.
# ---------------------------------------------------------------
$ws.Range("K11").Value = "I#"
$ws.Range("K11").Font.Bold = $false
$ws.Range("K11").HorizontalAlignment = -4108
$ws.Range("K11").VerticalAlignment = -4108
$ws.Range("K11").Borders.LineStyle = 1
$ws.Range("K11").Borders.Weight = 2
$ws.Range("K11").Borders.Color = 0

$ws.Range("L11").Value = "ITEM SET/STATE #"
$ws.Range("L11").HorizontalAlignment = -4131
$ws.Range("L11").Borders.LineStyle = 1
$ws.Range("L11").Borders.Weight = 2
$ws.Range("L11").Borders.Color = 0

# ---------------------------------------------------------------
# 6) Column layout: column A now has its own explicit width entry
#    separate from B:C (same visible width, ~10.83 chars).
# ---------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 10

# ---------------------------------------------------------------
# 7) Selection / active cell ends on A7, matching the saved view.
# ---------------------------------------------------------------
[void]$ws.Range("A7").Select()
